# Fixed up Myxicola and a few more loose ends in second review
#
# Materials sheet ("Materials"):
#  - Drop the Taxon_Local_ID column (it held ${iNaturalistTaxonId}, which is
#    no longer part of this sheet's mapping).
#  - Drop the suborder / infraorder / superfamily columns (their mapped
#    values ${suborder} / ${infraorder} / ${superfamily} are removed too).
#  - Rename the scientificNameAuthorship mapping from ${summary.Author} to
#    the correct ${summary.authority}.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$columnsToRemove = @("Taxon_Local_ID", "suborder", "infraorder", "superfamily")
foreach ($header in $columnsToRemove) {
    $cell = $ws.Rows(1).Find($header)
    if ($cell -ne $null) {
        $cell.EntireColumn.Delete()
    }
}

$authorityHeader = $ws.Rows(1).Find("scientificNameAuthorship")
if ($authorityHeader -ne $null) {
    $ws.Cells(2, $authorityHeader.Column).Value = '${summary.authority}'
}
